$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 26: "16x perspective" sample, only perspective-correcting every 16 pixels.
$ws.Range("A26").Value = "16x perspective"
$ws.Range("B26").Value = 243

# Extend the shared formulas down to row 26 (matches existing column formulas).
$ws.Range("C26").Formula = "=B26/30"
$ws.Range("D26").Formula = "=B26/`$B`$2"

# Update selection to reflect the new last-used cell.
$ws.Range("B26").Select()
